$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextValue $ws "D2" "60.183.49"
Set-TextValue $ws "E2" "  -6.17%  "
Set-TextValue $ws "D3" "3.300.87"
Set-TextValue $ws "E3" "  -5.09%  "
Set-TextValue $ws "E4" "  +0.08%  "
Set-TextValue $ws "D5" "562.78"
Set-TextValue $ws "E5" "  -3.81%  "
Set-TextValue $ws "D6" "129.19"
Set-TextValue $ws "E6" "  -1.98%  "
Set-TextValue $ws "E7" "  +0.05%  "
Set-TextValue $ws "D8" "3.300.03"
Set-TextValue $ws "E8" "  -5.10%  "
Set-TextValue $ws "D9" "0.471"
Set-TextValue $ws "E9" "  -2.40%  "
Set-TextValue $ws "D10" "7.42"
Set-TextValue $ws "E10" "  -3.72%  "
Set-TextValue $ws "D11" "0.116"
Set-TextValue $ws "E11" "  -5.24%  "
Set-TextValue $ws "D12" "0.371"
Set-TextValue $ws "E12" "  -3.88%  "
Set-TextValue $ws "D13" "3.867.55"
Set-TextValue $ws "E13" "  -5.02%  "
Set-TextValue $ws "E14" "  -0.25%  "
Set-TextValue $ws "D15" "3.299.74"
Set-TextValue $ws "E15" "  -5.10%  "
Set-TextValue $ws "E16" "  -6.06%  "
Set-TextValue $ws "D17" "60.496.92"
Set-TextValue $ws "E17" "  -5.64%  "
Set-TextValue $ws "D18" "24.11"
Set-TextValue $ws "E18" "  -4.15%  "
Set-TextValue $ws "D19" "5.62"
Set-TextValue $ws "E19" "  -1.05%  "
Set-TextValue $ws "D20" "13.23"
Set-TextValue $ws "E20" "  -1.24%  "
Set-TextValue $ws "D21" "8.92"
Set-TextValue $ws "E21" "  -10.64%  "
Set-TextValue $ws "D22" "349.14"
Set-TextValue $ws "E22" "  -9.29%  "
Set-TextValue $ws "D23" "0.553"
Set-TextValue $ws "E23" "  -2.47%  "
Set-TextValue $ws "E24" "  +0.02%  "
Set-TextValue $ws "D25" "3.432.37"
Set-TextValue $ws "E25" "  -5.11%  "
Set-TextValue $ws "D26" "68.93"
Set-TextValue $ws "E26" "  -7.55%  "
Set-TextValue $ws "D27" "0.0000107"
Set-TextValue $ws "E27" "  -3.83%  "
Set-TextValue $ws "D28" "1.00"
Set-TextValue $ws "E28" "  +0.23%  "
Set-TextValue $ws "D29" "7.30"
Set-TextValue $ws "E29" "  +2.74%  "
Set-TextValue $ws "E30" "  +1.40%  "
Set-TextValue $ws "D31" "7.81"
Set-TextValue $ws "E31" "  -1.93%  "
Set-TextValue $ws "E32" "  -2.00%  "
Set-TextValue $ws "E33" "  -5.61%  "
Set-TextValue $ws "E34" "  +0.00%  "
Set-TextValue $ws "D35" "3.327.40"
Set-TextValue $ws "E35" "  -5.06%  "
Set-TextValue $ws "D36" "22.56"
Set-TextValue $ws "E36" "  -1.76%  "
Set-TextValue $ws "D37" "5.33"
Set-TextValue $ws "D38" "6.75"
Set-TextValue $ws "E38" "  -0.20%  "
Set-TextValue $ws "D39" "1.47"
Set-TextValue $ws "E39" "  -1.54%  "
Set-TextValue $ws "D40" "156.99"
Set-TextValue $ws "E40" "  -3.13%  "
Set-TextValue $ws "D41" "0.0749"
Set-TextValue $ws "E41" "  -3.76%  "
Set-TextValue $ws "D42" "1.00"
Set-TextValue $ws "E42" "  +0.12%  "
Set-TextValue $ws "D43" "40.84"
Set-TextValue $ws "E43" "  -1.55%  "
Set-TextValue $ws "D44" "4.31"
Set-TextValue $ws "E44" "  -0.54%  "
Set-TextValue $ws "D45" "0.741"
Set-TextValue $ws "E45" "  -7.11%  "
Set-TextValue $ws "D46" "1.16"
Set-TextValue $ws "E46" "  +2.72%  "
Set-TextValue $ws "D47" "1.53"
Set-TextValue $ws "D48" "22.39"
Set-TextValue $ws "E48" "  -4.35%  "
Set-TextValue $ws "E49" "  -0.82%  "
Set-TextValue $ws "D50" "21.72"
Set-TextValue $ws "E50" "  +5.99%  "
Set-TextValue $ws "D51" "0.858"
Set-TextValue $ws "E51" "  -4.68%  "